$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Stato libero a seguito di scioglimento dell'unione"
$ws.Range("C8").Value = "Stato libero a seguito di decesso della parte unita civilmente"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8

$ws.Range("B10").Select()
